$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '66.807.70'
$ws.Cells.Item(2, 5).Value = '  -0.61%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.091.35'
$ws.Cells.Item(3, 5).Value = '  -0.75%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.01%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '576.07'
$ws.Cells.Item(5, 5).Value = '  -0.66%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '177.74'
$ws.Cells.Item(6, 5).Value = '  +2.46%  '

$ws.Cells.Item(7, 5).Value = '  -0.06%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '3.090.94'
$ws.Cells.Item(8, 5).Value = '  -0.64%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.513'
$ws.Cells.Item(9, 5).Value = '  -1.41%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '6.37'
$ws.Cells.Item(10, 5).Value = '  -2.37%  '

$ws.Cells.Item(11, 5).Value = '  -2.14%  '

$ws.Cells.Item(12, 5).Value = '  -2.82%  '

$ws.Cells.Item(13, 5).Value = '  -3.28%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '36.01'
$ws.Cells.Item(14, 5).Value = '  -2.48%  '

$ws.Cells.Item(15, 5).Value = '  -0.60%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.606.06'
$ws.Cells.Item(16, 5).Value = '  -0.70%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '66.739.41'
$ws.Cells.Item(17, 5).Value = '  -0.66%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '6.98'
$ws.Cells.Item(18, 5).Value = '  -1.82%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '16.70'
$ws.Cells.Item(19, 5).Value = '  +0.46%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.090.28'
$ws.Cells.Item(20, 5).Value = '  -0.69%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '481.09'
$ws.Cells.Item(21, 5).Value = '  -2.16%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '7.74'
$ws.Cells.Item(22, 5).Value = '  -2.28%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.689'
$ws.Cells.Item(23, 5).Value = '  -2.37%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '83.50'
$ws.Cells.Item(24, 5).Value = '  -0.54%  '

$ws.Cells.Item(25, 5).Value = '  -4.13%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.24'
$ws.Cells.Item(26, 5).Value = '  -2.78%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.13'
$ws.Cells.Item(27, 5).Value = '  -4.16%  '

$ws.Cells.Item(28, 5).Value = '  +0.02%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.94'
$ws.Cells.Item(29, 5).Value = '  -0.39%  '

$ws.Cells.Item(30, 5).Value = '  -4.21%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.60'
$ws.Cells.Item(31, 5).Value = '  -2.72%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '27.94'
$ws.Cells.Item(32, 5).Value = '  -1.67%  '

$ws.Cells.Item(33, 5).Value = '  -2.04%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0₃0941'
$ws.Cells.Item(34, 5).Value = '  -0.31%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  +0.04%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '48.33'
$ws.Cells.Item(36, 5).Value = '  +2.44%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.59'
$ws.Cells.Item(37, 5).Value = '  -4.88%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.940'
$ws.Cells.Item(38, 5).Value = '  -3.65%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.311'
$ws.Cells.Item(39, 5).Value = '  +0.27%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '48.99'
$ws.Cells.Item(40, 5).Value = '  -2.15%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.00'
$ws.Cells.Item(41, 5).Value = '  -2.28%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.123'
$ws.Cells.Item(42, 5).Value = '  -0.48%  '

$ws.Cells.Item(43, 5).Value = '  -2.08%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.67'
$ws.Cells.Item(44, 5).Value = '  +3.00%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.792.15'
$ws.Cells.Item(45, 5).Value = '  -0.61%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '371.25'
$ws.Cells.Item(46, 5).Value = '  -4.33%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '135.50'
$ws.Cells.Item(47, 5).Value = '  +0.26%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0343'
$ws.Cells.Item(48, 5).Value = '  -2.53%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '24.89'
$ws.Cells.Item(50, 5).Value = '  -0.66%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.24'
$ws.Cells.Item(51, 5).Value = '  +1.01%  '

# Reset D-column cells that were forced to text format back to the default style
$ws.Range("D2:D51").Style = "Normal"

